$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: __init__ / Attributes are set to parameter values. (date/management_fee values updated)
$ws.Range("F7").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2024, 5, 10), management_fee = 3"

# Row 8: __init__ / management fee has invalid type.
$ws.Range("F8").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2024, 3, 25), management_fee = 'three'"

# Row 9: get_service_charges / date created more than 10 years ago
$ws.Range("F9").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2014, 3, 25), management_fee = 3"

# Row 10: get_service_charges / date created exactly 10 years ago.
$ws.Range("F10").Value = "account_number = 350, client_number = 350, balance = 350, date_created = InvestmentAccount.TEN_YEARS_AGO, management_fee = 3"

# Row 11: get_service_charges / date created within last 10 years.
$ws.Range("F11").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2024, 3, 25), management_fee = 3"

# Row 12: __str__ / displays waived management fee when date created more than 10 years ago.
$ws.Range("F12").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2014, 3, 25), management_fee = 3"
$ws.Range("G12").Value = "Account Number: 350 Balance: `$350.00 Date Created: 2014-05-10 Management Fee: Waived Account Type: Investment"

# Row 13: __str__ / displays management fee when date created within last 10 years.
$ws.Range("F13").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2024, 3, 25), management_fee = 3"
$ws.Range("G13").Value = "Account Number: 350 Balance: `$350.00 Date Created: 2024-05-10 Management Fee: `$3.00 Account Type: Investment"

# Restore sheet view: scroll/selection
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("G12").Select()
